# Create data sources table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A / B first (rows 1-6), in natural fill order -----------------
$ws.Range("A1").Value = "Data"
$ws.Range("B1").Value = "Source"

$ws.Range("A2").Value = "Bathymetry"
$ws.Range("B2").Value = "poster"

$ws.Range("A3").Value = "Irradiance"

$ws.Range("A4").Value = "Lake Washington WQ"
$ws.Range("A5").Value = "Lake Washington meteorology"

$ws.Range("B4").Value = "King County buoy"
$ws.Range("B5").Value = "King County buoy"

# --- Citation header (column D) --------------------------------------------
$ws.Range("D1").Value = "Citation"

# --- Hyperlinks for the King County buoy rows (also sets cell text) --------
$ws.Hyperlinks.Add($ws.Range("C4"), "https://green2.kingcounty.gov/lake-buoy/Data.aspx")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://green2.kingcounty.gov/lake-buoy/Data.aspx")

# --- URL header (column C) --------------------------------------------------
$ws.Range("C1").Value = "URL"

# --- Smoke data row ----------------------------------------------------------
$ws.Range("A6").Value = "Smoke data"
$ws.Range("B6").Value = "NOAA"
$ws.Hyperlinks.Add($ws.Range("C6"), "https://satepsanone.nesdis.noaa.gov/pub/volcano/FIRE/HMS_ARCHIVE/")

# --- Column widths (best-fit sized to the longest entry in each column) ----
$ws.Columns("A:A").ColumnWidth = 27.66666666666667
$ws.Columns("C:C").ColumnWidth = 47.66666666666667

# --- Leave the selection where it lands after entering the table -----------
$ws.Range("A7").Select() | Out-Null
